# Generate Report for Handback
# Updates the timestamp strings recorded in the handback-status workbook.
#
# Mapping of changed shared strings -> worksheet cells (discovered via the
# tables / hyperlinks referencing each sheet's "datetime" columns):
#   Overview!G2  ("Latest HO Xliff Generate Date")      2016-08-27 01:03:31 -> 2016-08-27 01:04:27
#   zh-cn!H2     ("Correspond Handoff Datetime")         2016-08-27 01:03:27 -> 2016-08-27 01:04:22
#   zh-cn!K2     ("Correspond Handback DateTime")        2016-08-27 01:03:55 -> 2016-08-27 01:04:39
#   de-de!K2     ("Correspond Handback DateTime")        2016-08-27 01:04:03 -> 2016-08-27 01:04:45

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-27 01:04:27"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-27 01:04:22"
$wsZhCn.Range("K2").Value = "2016-08-27 01:04:39"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-08-27 01:04:45"
